$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column C: "label" -> "score"
$ws.Range("C1").Value = "score"

# Update row 2 sample data
$ws.Range("A2").Value = "This is a cat"
$ws.Range("B2").Value = "This is a dog"

# C2 must become a text value "0.7" (not a number) - force text formatting,
# assign, then reset the style so no leftover number-format style lingers.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0.7"
$ws.Range("C2").Style = "Normal"

# Remove the old sample rows 3-11, leaving just the header + one data row.
$ws.Range("A3:C11").EntireRow.Delete()
